$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns stay text-formatted (values like "4.40" or "0.850"
# must keep their trailing zeros, matching the original inline-string cells).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.175.33'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.642.61'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.07'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.514'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.58%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.09%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.98'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.14%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.874.08'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.648.91'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.14'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.33'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.195.01'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.29%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.86'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.84'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.31%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.40'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.65%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.72'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.54'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.19%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.74'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.19'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.262.02'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.32%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.850'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.09%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.809'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.37%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.783.96'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.69'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.75'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0106'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.40%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.63'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0974'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.09%  '
